$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.094.88"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "'3.152.17"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'592.26"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "'146.14"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'3.142.93"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "'5.88"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'3.674.28"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "'63.914.37"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'3.148.82"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'467.84"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'14.37"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'7.51"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "'2.34"
$ws.Range("E24").Value = "  +6.93%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.01"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "'81.32"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D28").Value = "'9.74"
$ws.Range("E28").Value = "  +8.32%  "
$ws.Range("E29").Value = "  +7.76%  "
$ws.Range("D30").Value = "'2.71"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'27.73"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "'0.0₃0839"
$ws.Range("E35").Value = "  -4.55%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "'6.16"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "'3.29"
$ws.Range("E39").Value = "  -4.67%  "
$ws.Range("D40").Value = "'463.75"
$ws.Range("D41").Value = "'51.40"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "'9.21"
$ws.Range("E42").Value = "  +5.37%  "
$ws.Range("E43").Value = "  +5.53%  "
$ws.Range("D44").Value = "'2.927.62"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "'0.0372"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "'40.24"
$ws.Range("E46").Value = "  +12.65%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'129.20"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("E51").Value = "  -0.73%  "
